# Move the table's surrounding text onto the same single-column
# "Title and Content" layout, stacking Title / body text / table
# vertically instead of the old two-column "Content with Caption"
# layout.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$m = $p.SlideMaster
$newLayout = $m.CustomLayouts.Item(2)   # "Title and Content"

# --- Preserve the table (it keeps its content/style, only moves). ---
$tbl = $s.Shapes.Item(3)
$tbl.Copy()
$tbl.Delete()

# --- Drop the old "half" body placeholder (idx=2) entirely. ---
$body = $s.Shapes.Item(2)
$body.Delete()
$body2 = $s.Shapes.Item(2)
$body2.Delete()

# --- Drop the old title placeholder so it can be re-synthesized with
#     a clean (layout-inherited) position. ---
$title = $s.Shapes.Item(1)
$title.Delete()
$title2 = $s.Shapes.Item(1)
$title2.Delete()

# --- Re-point the slide at "Title and Content"; this re-creates the
#     Title + Content placeholders (idx=1) with default positions. ---
$s.CustomLayout = $newLayout

# --- Title: re-enter the text, matching the original plain run. ---
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Hello"
$titleRange.IndentLevel = 0
$titleRange.ParagraphFormat.Bullet.Visible = 0

# --- Content placeholder (was "Text Placeholder 3"): text + position. ---
$bodyShape = $s.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "There"
$bodyRange.IndentLevel = 0
$bodyRange.ParagraphFormat.Bullet.Visible = 0
$bodyRange.ParagraphFormat.SpaceBefore = 30
$bodyRange.Font.Bold = 1

$bodyShape.Left = 36
$bodyShape.Top = 94
$bodyShape.Width = 648
$bodyShape.Height = 124

# --- Paste the table back, then move/resize it below the body text. ---
$pasted = $s.Shapes.Paste()
$tblShape = $pasted.Item(1)
$tblShape.Left = 36
$tblShape.Top = 228
$tblShape.Width = 648
$tblShape.Height = 124
$tblShape.Table.Columns.Item(1).Width = 648
